# Applies the "Updated cryptos list" data refresh to Sheet1.
# For each changed cell we set the new text value. Columns D/E in this
# sheet are always stored as text (prices use "." as a thousands
# separator and percentages keep surrounding padding spaces), so any
# D-column value that Excel would otherwise auto-recognize as a plain
# number is forced to stay text via NumberFormat "@" before assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.495.59'
$ws.Range("E2").Value = '  +0.88%  '
# Row 3
$ws.Range("D3").Value = '3.315.19'
$ws.Range("E3").Value = '  +0.14%  '
# Row 4
$ws.Range("E4").Value = '  +0.01%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.57'
$ws.Range("E5").Value = '  +2.49%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.65'
$ws.Range("E6").Value = '  -0.14%  '
# Row 7
$ws.Range("E7").Value = '  +6.00%  '
# Row 8
$ws.Range("E8").Value = '  +0.03%  '
# Row 9
$ws.Range("D9").Value = '3.313.24'
$ws.Range("E9").Value = '  +0.11%  '
# Row 10
$ws.Range("E10").Value = '  -0.42%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.83'
$ws.Range("E11").Value = '  +2.68%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.403'
$ws.Range("E12").Value = '  +0.45%  '
# Row 13
$ws.Range("D13").Value = '3.894.04'
$ws.Range("E13").Value = '  +0.17%  '
# Row 14
$ws.Range("E14").Value = '  -2.68%  '
# Row 15
$ws.Range("D15").Value = '66.469.21'
$ws.Range("E15").Value = '  +0.63%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.64'
$ws.Range("E16").Value = '  -0.08%  '
# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.336.11'
$ws.Range("E17").Value = '  +1.58%  '
# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000164'
$ws.Range("E18").Value = '  -0.98%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '424.73'
$ws.Range("E19").Value = '  -2.61%  '
# Row 20
$ws.Range("E20").Value = '  -2.86%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.47'
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.34'
$ws.Range("E22").Value = '  -2.81%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.77'
$ws.Range("E23").Value = '  -1.98%  '
# Row 24
$ws.Range("E24").Value = '  -0.03%  '
# Row 25
$ws.Range("E25").Value = '  -0.11%  '
# Row 26
$ws.Range("D26").Value = '3.467.86'
$ws.Range("E26").Value = '  +0.27%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.515'
$ws.Range("E27").Value = '  -0.81%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.204'
$ws.Range("E28").Value = '  +5.79%  '
# Row 29
$ws.Range("E29").Value = '  -0.56%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.08'
$ws.Range("E30").Value = '  +0.15%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.32%  '
# Row 32
$ws.Range("E32").Value = '  -1.58%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.40'
$ws.Range("E33").Value = '  -1.15%  '
# Row 34
$ws.Range("E34").Value = '  +0.03%  '
# Row 35
$ws.Range("E35").Value = '  -1.04%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.58'
$ws.Range("E36").Value = '  -2.49%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.19'
$ws.Range("E37").Value = '  -2.12%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '159.95'
$ws.Range("E38").Value = '  +0.17%  '
# Row 39
$ws.Range("E39").Value = '  -3.22%  '
# Row 40
$ws.Range("D40").Value = '2.860.20'
$ws.Range("E40").Value = '  +1.08%  '
# Row 41
$ws.Range("E41").Value = '  +0.90%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.36'
$ws.Range("E42").Value = '  -4.79%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.34'
$ws.Range("E43").Value = '  -1.92%  '
# Row 44
$ws.Range("E44").Value = '  -3.84%  '
# Row 45
$ws.Range("E45").Value = '  -1.45%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0658'
$ws.Range("E46").Value = '  -1.06%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.89'
$ws.Range("E47").Value = '  -4.61%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.30'
$ws.Range("E48").Value = '  -2.21%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.12'
$ws.Range("E49").Value = '  -4.10%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '311.71'
$ws.Range("E50").Value = '  -4.17%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0273'
$ws.Range("E51").Value = '  +0.56%  '
